$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their exact literal formatting (avoid numeric auto-conversion)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.882.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.16%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.496.07"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.16%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.53"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.44%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.53"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.84%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.36%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.10%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.82%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.04%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.62%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.37%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.938.60"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.50%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "58.815.73"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.20%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.03%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.93%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.501.51"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.74%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.01"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.48%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.24"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.19%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "322.83"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.84%  "

# Row 21
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.00%  "

# Row 22
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.93"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.13%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.01"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.56%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.420"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.56%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.83%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.57%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.51"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.10%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.02%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.30"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.22%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.56%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.42"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.46%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.16"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.19%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.13%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.33"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.84%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.18%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.05"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.50%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.97%  "

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.798"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.84%  "

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.56"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.60%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "280.84"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.06%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.41%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.25%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "129.51"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.78%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.89"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.52%  "

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.19%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0923"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.34%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0498"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.57%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.53%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.20"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.93%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.752.89"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.66%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.983"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.27%  "
